$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
${ws}.Range("E2").Value = 3
${ws}.Range("G2").Value = 15.829186
${ws}.Range("H2").Value = 47.487558
${ws}.Range("I2").Value = 0.01520167221269649
${ws}.Range("J2").Value = 0.01552195334947967
${ws}.Range("K2").Value = 3
${ws}.Range("M2").Value = 138.2190853333333
${ws}.Range("N2").Value = 414.657256
${ws}.Range("O2").Value = 0.2249223651785973
${ws}.Range("P2").Value = 0.2476599003709697
${ws}.Range("Q2").Value = 2187.895610491205
${ws}.Range("R2").Value = 19691.06049442085
${ws}.Range("S2").Value = 0.003419196068749456
${ws}.Range("T2").Value = 0.003844165420094975

# Row 3
${ws}.Range("E3").Value = 3
${ws}.Range("G3").Value = 15.829186
${ws}.Range("H3").Value = 47.487558
${ws}.Range("I3").Value = 0.01520167221269649
${ws}.Range("J3").Value = 0.01552195334947967
${ws}.Range("K3").Value = 3
${ws}.Range("M3").Value = 147.91433
${ws}.Range("N3").Value = 443.74299
${ws}.Range("O3").Value = 0.2406993279341593
${ws}.Range("P3").Value = 0.2650317656414439
${ws}.Range("Q3").Value = 2341.36344163538
${ws}.Range("R3").Value = 21072.27097471842
${ws}.Range("S3").Value = 0.00365903228507143
${ws}.Range("T3").Value = 0.004113810702416721

# Row 4
${ws}.Range("E4").Value = 3
${ws}.Range("G4").Value = 15.829186
${ws}.Range("H4").Value = 47.487558
${ws}.Range("I4").Value = 0.01520167221269649
${ws}.Range("J4").Value = 0.01552195334947967
${ws}.Range("K4").Value = 3
${ws}.Range("M4").Value = 74.27261733333333
${ws}.Range("N4").Value = 222.817852
${ws}.Range("O4").Value = 0.1208629960061633
${ws}.Range("P4").Value = 0.1330811078998542
${ws}.Range("Q4").Value = 1175.675074476157
${ws}.Range("R4").Value = 10581.07567028542
${ws}.Range("S4").Value = 0.00183731964793014
${ws}.Range("T4").Value = 0.002065678748518607

# Row 5
${ws}.Range("E5").Value = 3
${ws}.Range("G5").Value = 15.829186
${ws}.Range("H5").Value = 47.487558
${ws}.Range("I5").Value = 0.01520167221269649
${ws}.Range("J5").Value = 0.01552195334947967
${ws}.Range("K5").Value = 3
${ws}.Range("M5").Value = 84.85695366666668
${ws}.Range("N5").Value = 254.570861
${ws}.Range("O5").Value = 0.138086767645209
${ws}.Range("P5").Value = 0.1520460408212704
${ws}.Range("Q5").Value = 1343.216502983049
${ws}.Range("R5").Value = 12088.94852684744
${ws}.Range("S5").Value = 0.002099149778653251
${ws}.Range("T5").Value = 0.002360051552600842

# Row 6
${ws}.Range("E6").Value = 3
${ws}.Range("G6").Value = 15.829186
${ws}.Range("H6").Value = 47.487558
${ws}.Range("I6").Value = 0.01520167221269649
${ws}.Range("J6").Value = 0.01552195334947967
${ws}.Range("K6").Value = 2
${ws}.Range("M6").Value = 169.2560955
${ws}.Range("N6").Value = 338.512191
${ws}.Range("O6").Value = 0.275428543235871
${ws}.Range("P6").Value = 0.2021811852664618
${ws}.Range("Q6").Value = 2679.186217303263
${ws}.Range("R6").Value = 16075.11730381958
${ws}.Range("S6").Value = 0.004186974432292215
${ws}.Range("T6").Value = 0.003138246925848526

# Row 7
${ws}.Range("E7").Value = 3
${ws}.Range("G7").Value = 155.500389
${ws}.Range("H7").Value = 466.501167
${ws}.Range("I7").Value = 0.1493359129474374
${ws}.Range("J7").Value = 0.1524822428572096
${ws}.Range("K7").Value = 3
${ws}.Range("M7").Value = 138.2190853333333
${ws}.Range("N7").Value = 414.657256
${ws}.Range("O7").Value = 0.2249223651785973
${ws}.Range("P7").Value = 0.2476599003709697
${ws}.Range("Q7").Value = 21493.12153655753
${ws}.Range("R7").Value = 193438.0938290178
${ws}.Range("S7").Value = 0.03358898674624274
${ws}.Range("T7").Value = 0.03776373707435853

# Row 8
${ws}.Range("E8").Value = 3
${ws}.Range("G8").Value = 155.500389
${ws}.Range("H8").Value = 466.501167
${ws}.Range("I8").Value = 0.1493359129474374
${ws}.Range("J8").Value = 0.1524822428572096
${ws}.Range("K8").Value = 3
${ws}.Range("M8").Value = 147.91433
${ws}.Range("N8").Value = 443.74299
${ws}.Range("O8").Value = 0.2406993279341593
${ws}.Range("P8").Value = 0.2650317656414439
${ws}.Range("Q8").Value = 23000.73585367437
${ws}.Range("R8").Value = 207006.6226830693
${ws}.Range("S8").Value = 0.03594505388288231
${ws}.Range("T8").Value = 0.0404126380534137

# Row 9
${ws}.Range("E9").Value = 3
${ws}.Range("G9").Value = 155.500389
${ws}.Range("H9").Value = 466.501167
${ws}.Range("I9").Value = 0.1493359129474374
${ws}.Range("J9").Value = 0.1524822428572096
${ws}.Range("K9").Value = 3
${ws}.Range("M9").Value = 74.27261733333333
${ws}.Range("N9").Value = 222.817852
${ws}.Range("O9").Value = 0.1208629960061633
${ws}.Range("P9").Value = 0.1330811078998542
${ws}.Range("Q9").Value = 11549.42088738148
${ws}.Range("R9").Value = 103944.7879864333
${ws}.Range("S9").Value = 0.01804918585014289
${ws}.Range("T9").Value = 0.02029250581449208

# Row 10
${ws}.Range("E10").Value = 3
${ws}.Range("G10").Value = 155.500389
${ws}.Range("H10").Value = 466.501167
${ws}.Range("I10").Value = 0.1493359129474374
${ws}.Range("J10").Value = 0.1524822428572096
${ws}.Range("K10").Value = 3
${ws}.Range("M10").Value = 84.85695366666668
${ws}.Range("N10").Value = 254.570861
${ws}.Range("O10").Value = 0.138086767645209
${ws}.Range("P10").Value = 0.1520460408212704
${ws}.Range("Q10").Value = 13195.28930452165
${ws}.Range("R10").Value = 118757.6037406948
${ws}.Range("S10").Value = 0.02062131351225795
${ws}.Range("T10").Value = 0.02318432132198616

# Row 11
${ws}.Range("E11").Value = 3
${ws}.Range("G11").Value = 155.500389
${ws}.Range("H11").Value = 466.501167
${ws}.Range("I11").Value = 0.1493359129474374
${ws}.Range("J11").Value = 0.1524822428572096
${ws}.Range("K11").Value = 2
${ws}.Range("M11").Value = 169.2560955
${ws}.Range("N11").Value = 338.512191
${ws}.Range("O11").Value = 0.275428543235871
${ws}.Range("P11").Value = 0.2021811852664618
${ws}.Range("Q11").Value = 26319.38869087115
${ws}.Range("R11").Value = 157916.3321452269
${ws}.Range("S11").Value = 0.04113137295591154
${ws}.Range("T11").Value = 0.03082904059295911

# Row 12
${ws}.Range("E12").Value = 3
${ws}.Range("G12").Value = 179.1193723333333
${ws}.Range("H12").Value = 537.358117
${ws}.Range("I12").Value = 0.1720185728536685
${ws}.Range("J12").Value = 0.1756427994052303
${ws}.Range("K12").Value = 3
${ws}.Range("M12").Value = 138.2190853333333
${ws}.Range("N12").Value = 414.657256
${ws}.Range("O12").Value = 0.2249223651785973
${ws}.Range("P12").Value = 0.2476599003709697
${ws}.Range("Q12").Value = 24757.71580939411
${ws}.Range("R12").Value = 222819.442284547
${ws}.Range("S12").Value = 0.03869082426089399
${ws}.Range("T12").Value = 0.04349967820157755

# Row 13
${ws}.Range("E13").Value = 3
${ws}.Range("G13").Value = 179.1193723333333
${ws}.Range("H13").Value = 537.358117
${ws}.Range("I13").Value = 0.1720185728536685
${ws}.Range("J13").Value = 0.1756427994052303
${ws}.Range("K13").Value = 3
${ws}.Range("M13").Value = 147.91433
${ws}.Range("N13").Value = 443.74299
${ws}.Range("O13").Value = 0.2406993279341593
${ws}.Range("P13").Value = 0.2650317656414439
${ws}.Range("Q13").Value = 26494.32194870553
${ws}.Range("R13").Value = 238448.8975383498
${ws}.Range("S13").Value = 0.04140475487807124
${ws}.Range("T13").Value = 0.04655092124857413

# Row 14
${ws}.Range("E14").Value = 3
${ws}.Range("G14").Value = 179.1193723333333
${ws}.Range("H14").Value = 537.358117
${ws}.Range("I14").Value = 0.1720185728536685
${ws}.Range("J14").Value = 0.1756427994052303
${ws}.Range("K14").Value = 3
${ws}.Range("M14").Value = 74.27261733333333
${ws}.Range("N14").Value = 222.817852
${ws}.Range("O14").Value = 0.1208629960061633
${ws}.Range("P14").Value = 0.1330811078998542
${ws}.Range("Q14").Value = 13303.66459830052
${ws}.Range("R14").Value = 119732.9813847047
${ws}.Range("S14").Value = 0.02079068008379886
${ws}.Range("T14").Value = 0.0233747383394799

# Row 15
${ws}.Range("E15").Value = 3
${ws}.Range("G15").Value = 179.1193723333333
${ws}.Range("H15").Value = 537.358117
${ws}.Range("I15").Value = 0.1720185728536685
${ws}.Range("J15").Value = 0.1756427994052303
${ws}.Range("K15").Value = 3
${ws}.Range("M15").Value = 84.85695366666668
${ws}.Range("N15").Value = 254.570861
${ws}.Range("O15").Value = 0.138086767645209
${ws}.Range("P15").Value = 0.1520460408212704
${ws}.Range("Q15").Value = 15199.52427889208
${ws}.Range("R15").Value = 136795.7185100288
${ws}.Range("S15").Value = 0.02375348870030499
${ws}.Range("T15").Value = 0.02670579224832986

# Row 16
${ws}.Range("E16").Value = 3
${ws}.Range("G16").Value = 179.1193723333333
${ws}.Range("H16").Value = 537.358117
${ws}.Range("I16").Value = 0.1720185728536685
${ws}.Range("J16").Value = 0.1756427994052303
${ws}.Range("K16").Value = 2
${ws}.Range("M16").Value = 169.2560955
${ws}.Range("N16").Value = 338.512191
${ws}.Range("O16").Value = 0.275428543235871
${ws}.Range("P16").Value = 0.2021811852664618
${ws}.Range("Q16").Value = 30317.04558955073
${ws}.Range("R16").Value = 181902.2735373044
${ws}.Range("S16").Value = 0.04737882493059947
${ws}.Range("T16").Value = 0.03551166936726885

# Row 17
${ws}.Range("E17").Value = 3
${ws}.Range("G17").Value = 626.3728126666666
${ws}.Range("H17").Value = 1879.118438
${ws}.Range("I17").Value = 0.6015416194555684
${ws}.Range("J17").Value = 0.6142153852759307
${ws}.Range("K17").Value = 3
${ws}.Range("M17").Value = 138.2190853333333
${ws}.Range("N17").Value = 414.657256
${ws}.Range("O17").Value = 0.2249223651785973
${ws}.Range("P17").Value = 0.2476599003709697
${ws}.Range("Q17").Value = 86576.67724445401
${ws}.Range("R17").Value = 779190.0952000861
${ws}.Range("S17").Value = 0.1353001638013102
${ws}.Range("T17").Value = 0.1521165211237538

# Row 18
${ws}.Range("E18").Value = 3
${ws}.Range("G18").Value = 626.3728126666666
${ws}.Range("H18").Value = 1879.118438
${ws}.Range("I18").Value = 0.6015416194555684
${ws}.Range("J18").Value = 0.6142153852759307
${ws}.Range("K18").Value = 3
${ws}.Range("M18").Value = 147.91433
${ws}.Range("N18").Value = 443.74299
${ws}.Range("O18").Value = 0.2406993279341593
${ws}.Range("P18").Value = 0.2650317656414439
${ws}.Range("Q18").Value = 92649.5149158055
${ws}.Range("R18").Value = 833845.6342422495
${ws}.Range("S18").Value = 0.1447906635273811
${ws}.Range("T18").Value = 0.1627865880438196

# Row 19
${ws}.Range("E19").Value = 3
${ws}.Range("G19").Value = 626.3728126666666
${ws}.Range("H19").Value = 1879.118438
${ws}.Range("I19").Value = 0.6015416194555684
${ws}.Range("J19").Value = 0.6142153852759307
${ws}.Range("K19").Value = 3
${ws}.Range("M19").Value = 74.27261733333333
${ws}.Range("N19").Value = 222.817852
${ws}.Range("O19").Value = 0.1208629960061633
${ws}.Range("P19").Value = 0.1330811078998542
${ws}.Range("Q19").Value = 46522.34822319501
${ws}.Range("R19").Value = 418701.1340087552
${ws}.Range("S19").Value = 0.0727041223497994
${ws}.Range("T19").Value = 0.08174046396165664

# Row 20
${ws}.Range("E20").Value = 3
${ws}.Range("G20").Value = 626.3728126666666
${ws}.Range("H20").Value = 1879.118438
${ws}.Range("I20").Value = 0.6015416194555684
${ws}.Range("J20").Value = 0.6142153852759307
${ws}.Range("K20").Value = 3
${ws}.Range("M20").Value = 84.85695366666668
${ws}.Range("N20").Value = 254.570861
${ws}.Range("O20").Value = 0.138086767645209
${ws}.Range("P20").Value = 0.1520460408212704
${ws}.Range("Q20").Value = 53152.08874251502
${ws}.Range("R20").Value = 478368.7986826352
${ws}.Range("S20").Value = 0.08306493783468381
${ws}.Range("T20").Value = 0.0933890175427165

# Row 21
${ws}.Range("E21").Value = 3
${ws}.Range("G21").Value = 626.3728126666666
${ws}.Range("H21").Value = 1879.118438
${ws}.Range("I21").Value = 0.6015416194555684
${ws}.Range("J21").Value = 0.6142153852759307
${ws}.Range("K21").Value = 2
${ws}.Range("M21").Value = 169.2560955
${ws}.Range("N21").Value = 338.512191
${ws}.Range("O21").Value = 0.275428543235871
${ws}.Range("P21").Value = 0.2021811852664618
${ws}.Range("Q21").Value = 106017.4165993129
${ws}.Range("R21").Value = 636104.4995958777
${ws}.Range("S21").Value = 0.1656817319423939
${ws}.Range("T21").Value = 0.1241827946039841

# Row 22
${ws}.Range("E22").Value = 2
${ws}.Range("G22").Value = 64.4575005
${ws}.Range("H22").Value = 128.915001
${ws}.Range("I22").Value = 0.06190222253062919
${ws}.Range("J22").Value = 0.04213761911214986
${ws}.Range("K22").Value = 3
${ws}.Range("M22").Value = 138.2190853333333
${ws}.Range("N22").Value = 414.657256
${ws}.Range("O22").Value = 0.2249223651785973
${ws}.Range("P22").Value = 0.2476599003709697
${ws}.Range("Q22").Value = 8909.256761982875
${ws}.Range("R22").Value = 53455.54057189725
${ws}.Range("S22").Value = 0.01392319430140097
${ws}.Range("T22").Value = 0.0104357985511849

# Row 23
${ws}.Range("E23").Value = 2
${ws}.Range("G23").Value = 64.4575005
${ws}.Range("H23").Value = 128.915001
${ws}.Range("I23").Value = 0.06190222253062919
${ws}.Range("J23").Value = 0.04213761911214986
${ws}.Range("K23").Value = 3
${ws}.Range("M23").Value = 147.91433
${ws}.Range("N23").Value = 443.74299
${ws}.Range("O23").Value = 0.2406993279341593
${ws}.Range("P23").Value = 0.2650317656414439
${ws}.Range("Q23").Value = 9534.187999932163
${ws}.Range("R23").Value = 57205.12799959298
${ws}.Range("S23").Value = 0.01489982336075322
${ws}.Range("T23").Value = 0.01116780759321973

# Row 24
${ws}.Range("E24").Value = 2
${ws}.Range("G24").Value = 64.4575005
${ws}.Range("H24").Value = 128.915001
${ws}.Range("I24").Value = 0.06190222253062919
${ws}.Range("J24").Value = 0.04213761911214986
${ws}.Range("K24").Value = 3
${ws}.Range("M24").Value = 74.27261733333333
${ws}.Range("N24").Value = 222.817852
${ws}.Range("O24").Value = 0.1208629960061633
${ws}.Range("P24").Value = 0.1330811078998542
${ws}.Range("Q24").Value = 4787.427268899642
${ws}.Range("R24").Value = 28724.56361339785
${ws}.Range("S24").Value = 0.00748168807449207
${ws}.Range("T24").Value = 0.005607721035706974

# Row 25
${ws}.Range("E25").Value = 2
${ws}.Range("G25").Value = 64.4575005
${ws}.Range("H25").Value = 128.915001
${ws}.Range("I25").Value = 0.06190222253062919
${ws}.Range("J25").Value = 0.04213761911214986
${ws}.Range("K25").Value = 3
${ws}.Range("M25").Value = 84.85695366666668
${ws}.Range("N25").Value = 254.570861
${ws}.Range("O25").Value = 0.138086767645209
${ws}.Range("P25").Value = 0.1520460408212704
${ws}.Range("Q25").Value = 5469.667133397644
${ws}.Range("R25").Value = 32818.00280038586
${ws}.Range("S25").Value = 0.008547877819309015
${ws}.Range("T25").Value = 0.006406858155637083

# Row 26
${ws}.Range("E26").Value = 2
${ws}.Range("G26").Value = 64.4575005
${ws}.Range("H26").Value = 128.915001
${ws}.Range("I26").Value = 0.06190222253062919
${ws}.Range("J26").Value = 0.04213761911214986
${ws}.Range("K26").Value = 2
${ws}.Range("M26").Value = 169.2560955
${ws}.Range("N26").Value = 338.512191
${ws}.Range("O26").Value = 0.275428543235871
${ws}.Range("P26").Value = 0.2021811852664618
${ws}.Range("Q26").Value = 10909.8248603193
${ws}.Range("R26").Value = 43639.29944127719
${ws}.Range("S26").Value = 0.01704963897467391
${ws}.Range("T26").Value = 0.008519433776401172
